$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look numeric (e.g. "577.24") must be forced to
# stay text, matching the original inline-string cell type.
$textCells = @('D5','D6','D8','D13','D14','D19','D20','D21','D25','D26','D27','D31','D32','D38','D41','D46','D47','D49')
foreach ($addr in $textCells) { $ws.Range($addr).NumberFormat = "@" }

$ws.Range('D2').Value = '63.261.16'
$ws.Range('E2').Value = '  +1.21%  '
$ws.Range('D3').Value = '2.480.18'
$ws.Range('E3').Value = '  +3.05%  '
$ws.Range('E4').Value = '  -0.32%  '
$ws.Range('D5').Value = '577.24'
$ws.Range('E5').Value = '  +0.66%  '
$ws.Range('D6').Value = '146.62'
$ws.Range('E6').Value = '  +0.60%  '
$ws.Range('E7').Value = '  +0.21%  '
$ws.Range('D8').Value = '0.539'
$ws.Range('E8').Value = '  +0.12%  '
$ws.Range('D9').Value = '2.479.56'
$ws.Range('E9').Value = '  +1.92%  '
$ws.Range('E10').Value = '  +0.36%  '
$ws.Range('E11').Value = '  +1.72%  '
$ws.Range('E12').Value = '  +0.58%  '
$ws.Range('D13').Value = '0.354'
$ws.Range('E13').Value = '  +0.56%  '
$ws.Range('D14').Value = '28.57'
$ws.Range('E14').Value = '  +4.21%  '
$ws.Range('D16').Value = '2.931.10'
$ws.Range('E16').Value = '  +1.64%  '
$ws.Range('D17').Value = '63.178.49'
$ws.Range('E17').Value = '  +1.31%  '
$ws.Range('D18').Value = '2.480.12'
$ws.Range('E18').Value = '  +2.24%  '
$ws.Range('D19').Value = '8.16'
$ws.Range('E19').Value = '  +3.76%  '
$ws.Range('D20').Value = '11.03'
$ws.Range('E20').Value = '  +0.84%  '
$ws.Range('D21').Value = '330.35'
$ws.Range('E21').Value = '  +0.85%  '
$ws.Range('E22').Value = '  +8.88%  '
$ws.Range('E23').Value = '  +0.05%  '
$ws.Range('E24').Value = '  +0.22%  '
$ws.Range('D25').Value = '66.25'
$ws.Range('D26').Value = '9.85'
$ws.Range('E26').Value = '  +15.37%  '
$ws.Range('D27').Value = '661.55'
$ws.Range('E27').Value = '  +5.98%  '
$ws.Range('E28').Value = '  +1.71%  '
$ws.Range('D29').Value = '2.609.32'
$ws.Range('E29').Value = '  +1.88%  '
$ws.Range('E30').Value = '  +293.05%  '
$ws.Range('D31').Value = '1.49'
$ws.Range('E31').Value = '  +5.50%  '
$ws.Range('D32').Value = '8.10'
$ws.Range('E32').Value = '  -1.04%  '
$ws.Range('E33').Value = '  +1.06%  '
$ws.Range('E34').Value = '  -3.44%  '
$ws.Range('E35').Value = '  +4.22%  '
$ws.Range('E36').Value = '  +0.31%  '
$ws.Range('E37').Value = '  +0.87%  '
$ws.Range('D38').Value = '5.47'
$ws.Range('E38').Value = '  +1.71%  '
$ws.Range('E39').Value = '  -0.44%  '
$ws.Range('E40').Value = '  +0.98%  '
$ws.Range('D41').Value = '150.39'
$ws.Range('E41').Value = '  -0.97%  '
$ws.Range('E42').Value = '  -2.41%  '
$ws.Range('E43').Value = '  +0.60%  '
$ws.Range('E44').Value = '  +0.01%  '
$ws.Range('D45').Value = '0.0₆0312'
$ws.Range('E45').Value = '  -46.95%  '
$ws.Range('D46').Value = '156.08'
$ws.Range('E46').Value = '  +8.21%  '
$ws.Range('D47').Value = '15.25'
$ws.Range('E47').Value = '  +3.41%  '
$ws.Range('E48').Value = '  +0.74%  '
$ws.Range('D49').Value = '20.48'
$ws.Range('E49').Value = '  -0.10%  '
$ws.Range('E50').Value = '  +1.88%  '
$ws.Range('E51').Value = '  +0.13%  '

# Restore default (style-less) formatting now that the values are locked in as text.
foreach ($addr in $textCells) { $ws.Range($addr).Style = "Normal" }
